$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ranges = @(
    "G2:K2",
    "I3:K3",
    "G4:K4",
    "I5:K5",
    "G6:K6",
    "I7:K7",
    "G8:K8",
    "I9:K9",
    "G10:K10",
    "I11:K11",
    "G12:K12",
    "I13:K13",
    "G14:K14",
    "I15:K15",
    "G16:K16",
    "I17:K17",
    "K18:K18",
    "I19:K19",
    "K20:K20",
    "I21:K21",
    "K22:K22",
    "J23:K23",
    "I24:K24",
    "K26:K26",
    "J27:K27",
    "I28:K28",
    "K30:K30",
    "J31:K31",
    "I32:K32",
    "K34:K34",
    "J35:K35",
    "I36:K36",
    "K38:K38",
    "J39:K39",
    "I40:K40",
    "K42:K42",
    "J43:K43",
    "I44:J44"
)

foreach ($r in $ranges) {
    $ws.Range($r).ClearContents()
}
